$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $n = $parts.Count
        if ($n -gt 1) {
            for ($i = 0; $i -lt $n; $i++) {
                for ($j = 0; $j -lt ($n - $i - 1); $j++) {
                    $a = $parts[$j]
                    $b = $parts[$j + 1]
                    if ($a.CompareTo($b) -gt 0) {
                        $parts[$j] = $b
                        $parts[$j + 1] = $a
                    }
                }
            }
            $newVal = [string]::Join(", ", $parts)
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
